$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-21"

# Update the header label cell (A1) to reflect the new "through" date
$ws.Range("I1").Value = "2022 (through 11-21)"

# Update the data values for November (row 12) and Total (row 14)
$ws.Range("I12").Value = 81
$ws.Range("I14").Value = 1478
